$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "290×8="
$t.Cell(1,2).Range.Text = "147×7="
$t.Cell(1,3).Range.Text = "993×9="
$t.Cell(1,4).Range.Text = "789×7="
$t.Cell(1,5).Range.Text = "596×9="

# Row 5
$t.Cell(5,1).Range.Text = "221×3="
$t.Cell(5,2).Range.Text = "944×9="
$t.Cell(5,3).Range.Text = "573×3="
$t.Cell(5,4).Range.Text = "138×3="
$t.Cell(5,5).Range.Text = "336×5="

# Row 10
$t.Cell(10,1).Range.Text = "563×9="
$t.Cell(10,2).Range.Text = "990×2="
$t.Cell(10,3).Range.Text = "867×9="
$t.Cell(10,4).Range.Text = "289×8="
$t.Cell(10,5).Range.Text = "812×3="

# Row 15
$t.Cell(15,1).Range.Text = "520×2="
$t.Cell(15,2).Range.Text = "678×4="
$t.Cell(15,3).Range.Text = "479×9="
$t.Cell(15,4).Range.Text = "709×9="
$t.Cell(15,5).Range.Text = "112×9="

# Row 20
$t.Cell(20,1).Range.Text = "162×4="
$t.Cell(20,2).Range.Text = "203×9="
$t.Cell(20,3).Range.Text = "168×8="
$t.Cell(20,4).Range.Text = "699×4="
$t.Cell(20,5).Range.Text = "538×4="
